$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 42; existing rows 42-167 shift down to 43-168
$ws.Rows.Item(42).Insert()

# Populate the new row 42 with the new data record
$ws.Cells.Item(42, 1).Value = 8
$ws.Cells.Item(42, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(42, 3).Value = "Coquimbo"
$ws.Cells.Item(42, 4).Value = 44620
$ws.Cells.Item(42, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(42, 5).Value = 4
$ws.Cells.Item(42, 6).Value = 100112037
$ws.Cells.Item(42, 7).Value = "Cebollín"
$ws.Cells.Item(42, 8).Value = "Sin especificar"
$ws.Cells.Item(42, 9).Value = "Primera"
$ws.Cells.Item(42, 10).Value = 1000
$ws.Cells.Item(42, 11).Value = 1000
$ws.Cells.Item(42, 12).Value = 1200
$ws.Cells.Item(42, 13).Value = 1100
$ws.Cells.Item(42, 14).Value = "`$/paquete 6 unidades"
$ws.Cells.Item(42, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(42, 16).Value = 183
$ws.Cells.Item(42, 17).Value = 6
$ws.Cells.Item(42, 18).Value = "Hortaliza"
